$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the last sheet, named "Rozetka checklist" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Rozetka checklist"

# --- Checklist content (column A = section headers, column B = checklist items) ---

$rows = @(
    @('Главная страница', ''),
    @('', 'Сайт быстро загружается'),
    @('', 'Главное меню доступно и легко обнаруживается'),
    @('', 'Отображается логотип и брендинг'),
    @('', 'Баннеры и акции на главной странице привлекают внимание'),
    @('', 'Поиск находится на видном месте и легко доступен'),
    @('Страница категории товаров', ''),
    @('', 'Страница быстро загружается'),
    @('', 'Главное меню отображает категорию товаров'),
    @('', 'Категории товаров ясно указаны в главном меню'),
    @('', 'Изображения товаров ясно отображают товар и его детали'),
    @('', 'Фильтры на странице позволяют уточнить результаты поиска'),
    @('Страница товара', ''),
    @('', 'Страница товара быстро загружается'),
    @('', 'Отображается заголовок и цена товара'),
    @('', 'Изображения товара ясно отображают товар и его детали'),
    @('', 'Описание товара содержит достаточно информации'),
    @('', 'Кнопка "Купить" ясно отображается на странице товара'),
    @('Страница корзины', ''),
    @('', 'Страница корзины быстро загружается'),
    @('', 'Все добавленные товары отображаются в корзине'),
    @('', 'Общая стоимость заказа и количество товаров ясно указаны'),
    @('', 'Кнопка "Оформить заказ" ясно отображается на странице корзины'),
    @('Страница оформления заказа', ''),
    @('', 'Страница оформления заказа быстро загружается'),
    @('', 'Пользователь может легко вводить свои данные'),
    @('', 'Система оплаты безопасна и надежна'),
    @('', 'Политика доставки и возврата легко доступна и понятна'),
    @('', 'Общая стоимость заказа и количество товаров ясно указаны'),
)

$r = 1
foreach ($pair in $rows) {
    $a = $pair[0]
    $b = $pair[1]
    if ($a -ne "") {
        $ws.Cells.Item($r, 1).Value = $a
    }
    if ($b -ne "") {
        $ws.Cells.Item($r, 2).Value = $b
    }
    $r = $r + 1
}

# --- Formatting: apply to the full used range A1:B29 (also touches blank cells) ---
$rng = $ws.Range("A1:B29")
$rng.Font.Name = "Times New Roman"
$rng.Font.Size = 14
$rng.Interior.Pattern = 1
$rng.Interior.ThemeColor = 2
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4108
$rng.ShrinkToFit = $true
$rng.NumberFormat = "@"

# --- Column widths (engine quantizes to 1/6 character units; closest achievable) ---
$ws.Columns.Item(1).ColumnWidth = 26.5
$ws.Columns.Item(2).ColumnWidth = 89

# --- View state: scroll/selection and make this the active sheet/tab ---
$ws.Range("C29").Select()
$ws.Activate()
